$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14 (old rows 14-15 shift down to 15-16)
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the latest weekly price data
$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 45154
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100114002
$ws.Range("G14").Value = "Camote"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 16500
$ws.Range("L14").Value = 17000
$ws.Range("M14").Value = 16750
$ws.Range("N14").Value = '$/malla 18 kilos'
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 931
$ws.Range("Q14").Value = 18
$ws.Range("R14").Value = "Hortaliza"
